$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2230
$ws.Range("I32").Value = 2331.6
$ws.Range("J32").Value = 2117.111
$ws.Range("K32").Value = 2331.6
$ws.Range("L32").Value = 2117.111
$ws.Range("M32").Value = -2005.6
$ws.Range("N32").Value = -2769.111
$ws.Range("H51").Value = 3117.647
$ws.Range("I51").Value = 2750
$ws.Range("J51").Value = 3230.7693
$ws.Range("K51").Value = 2750
$ws.Range("L51").Value = 3230.7693
$ws.Range("M51").Value = -2266
$ws.Range("N51").Value = -4198.7693
$ws.Range("H125").Value = 1731.5555
$ws.Range("I125").Value = 1292.8
$ws.Range("J125").Value = 2280
$ws.Range("K125").Value = 11635.2
$ws.Range("L125").Value = 20520
$ws.Range("M125").Value = -9175.199999999999
$ws.Range("N125").Value = -25440
$ws.Range("H132").Value = 273370.72
$ws.Range("I132").Value = 288907.9
$ws.Range("J132").Value = 1470
$ws.Range("K132").Value = 866723.7000000001
$ws.Range("L132").Value = 4410
$ws.Range("M132").Value = -864193.7000000001
$ws.Range("N132").Value = -9470
$ws.Range("H135").Value = 1713.5625
$ws.Range("I135").Value = 643.03705
$ws.Range("J135").Value = 7494.4
$ws.Range("K135").Value = 5787.33345
$ws.Range("L135").Value = 67449.59999999999
$ws.Range("M135").Value = -3252.33345
$ws.Range("N135").Value = -72519.59999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88
$ws.Range("H32").Value = 8648.536
$ws.Range("I32").Value = 4851.0527
$ws.Range("J32").Value = 26686.584
$ws.Range("K32").Value = 4851.0527
$ws.Range("L32").Value = 26686.584
$ws.Range("M32").Value = -4564.0527
$ws.Range("N32").Value = -27260.584
$ws.Range("H61").Value = 2464.718
$ws.Range("J61").Value = 2941.0833
$ws.Range("L61").Value = 2941.0833
$ws.Range("N61").Value = -3365.0833
$ws.Range("H74").Value = 2503352.2
$ws.Range("I74").Value = 3448806.8
$ws.Range("K74").Value = 3448806.8
$ws.Range("M74").Value = -3447932.8
$ws.Range("H77").Value = 2503352.2
$ws.Range("I77").Value = 3448806.8
$ws.Range("K77").Value = 17244034
$ws.Range("M77").Value = -17239666
$ws.Range("H95").Value = 34062.668
$ws.Range("J95").Value = 34062.668
$ws.Range("L95").Value = 34062.668
$ws.Range("N95").Value = -39554.668
$ws.Range("H102").Value = 2270
$ws.Range("I102").Value = 2203.6365
$ws.Range("K102").Value = 2203.6365
$ws.Range("M102").Value = -581.6365000000001
$ws.Range("H132").Value = 10239883
$ws.Range("I132").Value = 13783523
$ws.Range("J132").Value = 2699.6667
$ws.Range("K132").Value = 41350569
$ws.Range("L132").Value = 8099.000100000001
$ws.Range("M132").Value = -41348039
$ws.Range("N132").Value = -13159.0001
$ws.Range("H136").Value = 2464.718
$ws.Range("J136").Value = 2941.0833
$ws.Range("L136").Value = 8823.249899999999
$ws.Range("N136").Value = -13923.2499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 1000
$ws.Range("K10").Value = 1000
$ws.Range("M10").Value = -860
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H56").Value = 49000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 49000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 49000
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -50478
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H103").Value = 23555.666
$ws.Range("J103").Value = 23555.666
$ws.Range("L103").Value = 23555.666
$ws.Range("N103").Value = -25899.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1912.3658
$ws.Range("I134").Value = 1930.3334
$ws.Range("J134").Value = 1863.3636
$ws.Range("K134").Value = 5791.0002
$ws.Range("L134").Value = 5590.0908
$ws.Range("M134").Value = -3256.0002
$ws.Range("N134").Value = -10660.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 842.54285
$ws.Range("I68").Value = 599.2653
$ws.Range("J68").Value = 1410.1904
$ws.Range("K68").Value = 1797.7959
$ws.Range("L68").Value = 4230.5712
$ws.Range("M68").Value = -986.7959000000001
$ws.Range("N68").Value = -5852.5712
$ws.Range("H71").Value = 842.54285
$ws.Range("I71").Value = 599.2653
$ws.Range("J71").Value = 1410.1904
$ws.Range("K71").Value = 5393.3877
$ws.Range("L71").Value = 12691.7136
$ws.Range("M71").Value = -1337.3877
$ws.Range("N71").Value = -20803.7136
$ws.Range("H97").Value = 2106.5557
$ws.Range("I97").Value = 151.5
$ws.Range("J97").Value = 2665.1428
$ws.Range("K97").Value = 454.5
$ws.Range("L97").Value = 7995.428400000001
$ws.Range("M97").Value = 41.5
$ws.Range("N97").Value = -8987.428400000001
$ws.Range("H107").Value = 28640.07
$ws.Range("I107").Value = 25594.95
$ws.Range("J107").Value = 32446.469
$ws.Range("K107").Value = 76784.85000000001
$ws.Range("L107").Value = 97339.40700000001
$ws.Range("M107").Value = -74864.85000000001
$ws.Range("N107").Value = -101179.407
$ws.Range("H131").Value = 1112836.4
$ws.Range("I131").Value = 881.1177
$ws.Range("J131").Value = 1371784.9
$ws.Range("K131").Value = 2643.3531
$ws.Range("L131").Value = 4115354.7
$ws.Range("M131").Value = 2396.6469
$ws.Range("N131").Value = -4125434.7
$ws.Range("H136").Value = 2093.1538
$ws.Range("I136").Value = 1522.375
$ws.Range("J136").Value = 3006.4
$ws.Range("K136").Value = 4567.125
$ws.Range("L136").Value = 9019.200000000001
$ws.Range("M136").Value = 532.875
$ws.Range("N136").Value = -19219.2
$ws.Range("H139").Value = 18521508
$ws.Range("I139").Value = 33335048
$ws.Range("J139").Value = 4583.3335
$ws.Range("K139").Value = 100005144
$ws.Range("L139").Value = 13750.0005
$ws.Range("M139").Value = -100000004
$ws.Range("N139").Value = -24030.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 49564
$ws.Range("J104").Value = 49564
$ws.Range("L104").Value = 49564
$ws.Range("N104").Value = -56552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2399.818
$ws.Range("I7").Value = 2399.818
$ws.Range("K7").Value = 2399.818
$ws.Range("M7").Value = -2287.818
$ws.Range("H22").Value = 386
$ws.Range("I22").Value = 299.77777
$ws.Range("J22").Value = 580
$ws.Range("K22").Value = 299.77777
$ws.Range("L22").Value = 580
$ws.Range("M22").Value = -4.777769999999975
$ws.Range("N22").Value = -1170
$ws.Range("H27").Value = 386
$ws.Range("I27").Value = 299.77777
$ws.Range("J27").Value = 580
$ws.Range("K27").Value = 299.77777
$ws.Range("L27").Value = 580
$ws.Range("M27").Value = -192.77777
$ws.Range("N27").Value = -794
$ws.Range("H97").Value = 13000
$ws.Range("J97").Value = 13000
$ws.Range("L97").Value = 13000
$ws.Range("N97").Value = -14982
$ws.Range("H126").Value = 2399.818
$ws.Range("I126").Value = 2399.818
$ws.Range("K126").Value = 7199.454000000001
$ws.Range("M126").Value = -4729.454000000001
$ws.Range("H132").Value = 4485.7466
$ws.Range("I132").Value = 4896.9316
$ws.Range("J132").Value = 3815.6667
$ws.Range("K132").Value = 14690.7948
$ws.Range("L132").Value = 11447.0001
$ws.Range("M132").Value = -12160.7948
$ws.Range("N132").Value = -16507.0001
$ws.Range("H136").Value = 2071.4482
$ws.Range("I136").Value = 1503.125
$ws.Range("J136").Value = 4799.4
$ws.Range("K136").Value = 4509.375
$ws.Range("L136").Value = 14398.2
$ws.Range("M136").Value = -1959.375
$ws.Range("N136").Value = -19498.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H104").Value = 33317.25
$ws.Range("J104").Value = 33317.25
$ws.Range("L104").Value = 33317.25
$ws.Range("N104").Value = -40305.25
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = 0
